# Replace backslashes with forward slashes in the video path cells
# on the "test_schedule_1" worksheet, and move the active selection
# on that sheet from D13 to C4 (its last used cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_schedule_1")

$ws.Range("C2").Value = "test_files/vids/test.mp4"
$ws.Range("C3").Value = "test_files/vids2/test2.mp4"
$ws.Range("C4").Value = "test_files/vids/test4.mp4"

$ws.Activate()
$ws.Range("C4").Select()
